# Apply hybrid bold + color (2C3E50) highlighting to quantitative impact
# metrics (percentages, dollar amounts, large numbers) across the
# achievements / work-experience bullet paragraphs, matching the target
# diff exactly. Each target paragraph currently holds a single plain run;
# locating a metric substring with a paragraph-scoped Find and flipping
# its Font Bold/Color causes the run to split into plain / highlighted
# runs automatically (mirroring Word's own Find-and-format behaviour), so
# no manual run-splitting/InsertAfter bookkeeping is required.

$d = $word.ActiveDocument

# OOXML stores <w:color w:val="2C3E50"/> (plain RGB hex), but the COM
# Font.Color property takes a BGR-packed long (0x00BBGGRR) - the same
# encoding VBA's RGB(r,g,b) produces - so the bytes of the target hex
# must be re-packed in reverse order before assigning.
function Get-BgrColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$metricColor = Get-BgrColor("2C3E50")

# Locate the single paragraph whose text contains $anchor (a substring
# unique to that bullet) - more resilient than a hard-coded paragraph
# index.
function Find-Paragraph([string]$anchor) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text.Contains($anchor)) {
            return $para
        }
    }
    return $null
}

# Within $paragraph, find the (unique) occurrence of $metricText and make
# it bold + the highlight color.
function Set-MetricBold($paragraph, [string]$metricText) {
    $r = $paragraph.Range.Duplicate
    $ok = $r.Find.Execute($metricText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $r.Font.Bold = 1
        $r.Font.Color = $metricColor
    }
    return $ok
}

# 1) "• Discovered systematic race coding errors ... from 23% to 64%"
#    (anchor is the distinctive "developed geospatial machine learning
#    algorithms..." clause so it cannot match either the shorter
#    "Discovered systematic race coding errors affecting all Black and
#    Asian-American voters" bullet under KEY ACHIEVEMENTS AND IMPACT, or
#    the unrelated "Machine learning platform that discovered systematic
#    coding errors and improved demographic classification accuracy from
#    23% to 64%" paragraph under KEY PROJECTS)
$p = Find-Paragraph("developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%")
Set-MetricBold $p "23%" | Out-Null
Set-MetricBold $p "64%" | Out-Null

# 2) "• Achieved 87% prediction accuracy ... ±4.2% to ±2.1%"
$p = Find-Paragraph("reducing polling error margins")
Set-MetricBold $p "87%" | Out-Null
Set-MetricBold $p "71%" | Out-Null
Set-MetricBold $p "±4.2%" | Out-Null
Set-MetricBold $p "±2.1%" | Out-Null

# 3) "• Wrote RFP and analyzed bids from 1,200 vendors ..."
$p = Find-Paragraph("Wrote RFP and analyzed bids")
Set-MetricBold $p "1,200" | Out-Null

# 4) "... became the $400M Polling Consortium Database ... now valued at $1B+"
$p = Find-Paragraph("Polling Consortium Database")
Set-MetricBold $p '$400M' | Out-Null
Set-MetricBold $p '$1B' | Out-Null

# 5) "• Algorithm reduced mapping costs by 73.5%, saving ... $4.7M" (KEY ACHIEVEMENTS bullet)
$p = Find-Paragraph("Algorithm reduced mapping costs")
Set-MetricBold $p "73.5%" | Out-Null
Set-MetricBold $p '$4.7M' | Out-Null

# 6) "• Achieved 87% prediction accuracy ... industry standard of 71%" (KEY ACHIEVEMENTS bullet,
#    shorter variant without the "reducing polling error margins" clause)
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t.Contains("Achieved") -and $t.Contains("87%") -and -not $t.Contains("reducing polling error margins")) {
        Set-MetricBold $para "87%" | Out-Null
        Set-MetricBold $para "71%" | Out-Null
    }
}
